$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (prevents Excel from auto-converting
# number-like strings such as "6.45" or "513.60" into numeric cells).
# Uses a scratch cell formatted as Text, then copies only the value
# across so the destination cell keeps its original (default) style.
function Set-TextValue {
    param($CellRef, $Text)
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $Text
    $scratch.Copy()
    $ws.Range($CellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

Set-TextValue "D2" "57.758.23"
$ws.Range("E2").Value = "  -2.91%  "
Set-TextValue "D3" "2.555.01"
$ws.Range("E3").Value = "  -3.59%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue "D5" "513.60"
$ws.Range("E5").Value = "  -3.23%  "
Set-TextValue "D6" "137.52"
$ws.Range("E6").Value = "  -5.91%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -2.59%  "
Set-TextValue "D9" "2.572.21"
$ws.Range("E9").Value = "  -3.59%  "
Set-TextValue "D10" "6.45"
$ws.Range("E10").Value = "  -3.78%  "
Set-TextValue "D11" "0.0980"
$ws.Range("E11").Value = "  -6.22%  "
$ws.Range("E12").Value = "  -4.57%  "
$ws.Range("E13").Value = "  +0.61%  "
Set-TextValue "D14" "3.002.42"
$ws.Range("E14").Value = "  -3.74%  "
Set-TextValue "D15" "57.727.15"
$ws.Range("E15").Value = "  -2.92%  "
Set-TextValue "D16" "19.89"
$ws.Range("E16").Value = "  -4.85%  "
Set-TextValue "D17" "2.555.70"
$ws.Range("E17").Value = "  -4.45%  "
$ws.Range("E18").Value = "  -5.29%  "
Set-TextValue "D19" "330.96"
$ws.Range("E19").Value = "  -3.55%  "
Set-TextValue "D20" "4.24"
$ws.Range("E20").Value = "  -4.86%  "
Set-TextValue "D21" "9.98"
$ws.Range("E21").Value = "  -6.29%  "
Set-TextValue "D22" "6.28"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("E23").Value = "  -0.07%  "
Set-TextValue "D24" "65.34"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  -5.76%  "
$ws.Range("E28").Value = "  -5.07%  "
Set-TextValue "D29" "0.998"
$ws.Range("E29").Value = "  -0.01%  "
Set-TextValue "D30" "0.0₃0698"
$ws.Range("E30").Value = "  -13.20%  "
Set-TextValue "D31" "5.86"
$ws.Range("E31").Value = "  -8.27%  "
Set-TextValue "D32" "148.87"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("E33").Value = "  -4.45%  "
Set-TextValue "D34" "18.43"
$ws.Range("E34").Value = "  -3.37%  "
$ws.Range("E35").Value = "  -8.03%  "
$ws.Range("E36").Value = "  -7.94%  "
Set-TextValue "D37" "35.95"
$ws.Range("E37").Value = "  -1.69%  "
Set-TextValue "D38" "0.820"
$ws.Range("E38").Value = "  -5.08%  "
Set-TextValue "D39" "0.814"
$ws.Range("E39").Value = "  -5.76%  "
Set-TextValue "D40" "1.39"
$ws.Range("E40").Value = "  -6.59%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D41" "0.998"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D42" "3.45"
$ws.Range("E42").Value = "  -5.51%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D43" "10.70"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D44" "266.69"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D45" "0.581"
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D46" "0.0934"
$ws.Range("E46").Value = "  -4.69%  "
$ws.Range("E47").Value = "  -5.16%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "18.10"
$ws.Range("E48").Value = "  -7.01%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D49" "1.944.59"
$ws.Range("E49").Value = "  -4.94%  "
Set-TextValue "D50" "0.0214"
$ws.Range("E50").Value = "  -6.52%  "
$ws.Range("E51").Value = "  -9.00%  "

$excel.CutCopyMode = $false

